$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Resize the table to include the new (9th) column first
$tbl = $ws.ListObjects.Item("Tabla1")
$tbl.Resize($ws.Range("A1:I3"))

# Add data for new column first so "datapro" gets the lower shared-string index
$ws.Range("I2").Value = "datapro"
$ws.Range("I3").Value = "datapro"

# Add new header for the nombre_proveedor column (also syncs the table column name)
$ws.Range("I1").Value = "nombre_proveedor"

# Set column width for column I
$ws.Range("I1").EntireColumn.ColumnWidth = 24.7109375

# Update selection
$ws.Range("A1:I3").Select()
